$wb = $excel.ActiveWorkbook

# --- Sheet 1: covariate_importance ---
$ws1 = $wb.Worksheets.Item("covariate_importance")

$sheet1Data = @(
    @("state", 100, 100, 100),
    @("region", 95, 85, 85),
    @("rplthemes", 75, 61, 61),
    @("derivedtotalenrolled", 63, 44, 44),
    @("percentblackorafricanamerican", 62, 42, 42),
    @("percentstudentsfreereducedlunch", 62, 41, 41),
    @("percentamericanindianoralaskanative", 50, 20, 20),
    @("percenthispaniclatino", 43, 20, 20),
    @("percenttwoormoreraces", 37, 17, 17),
    @("cntycaseschange", 43, 15, 14),
    @("locale", 27, 14, 14),
    @("schoollevel", 28, 12, 12),
    @("percentasian", 30, 11, 11),
    @("percentwhite", 39, 9, 9),
    @("percentnativehawaiianorotherpacificislander", 22, 4, 4),
    @("percentnotspecified", 8, 5, 0)
)

for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $row = $i + 2
    $rec = $sheet1Data[$i]
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $ws1.Cells.Item($row, 4).Value = $rec[3]
}

# --- Sheet 2: strategy_importance ---
$ws2 = $wb.Worksheets.Item("strategy_importance")

$sheet2Data = @(
    @("hvacsystems", 100, 100, 100),
    @("contacttracing", 78, 79, 76),
    @("masks", 44, 48, 43),
    @("vaccination", 40, 37, 39),
    @("screeningtestingforstudents", 38, 38, 38),
    @("cleaning", 39, 35, 35),
    @("quarantine", 24, 20, 22),
    @("physicaldistancing", 23, 21, 18),
    @("hepafilters", 6, 22, 4),
    @("stayhome", 0, 0, 0)
)

for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $row = $i + 2
    $rec = $sheet2Data[$i]
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
}
